$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new PR log row (row 3) with data for PR #23
$ws.Range("A3").Value = 23
$ws.Range("B3").Value = "Update index.py"
$ws.Range("C3").Value = "riya-morankar"
$ws.Range("D3").Value = "N/A"
$ws.Range("E3").Value = "edit1 to main"

# Use a leading apostrophe so Excel stores the date-looking text as a
# literal string instead of auto-converting it to a date serial number,
# then reset the cell style so the quote-prefix formatting doesn't stick.
$ws.Range("F3").Value = "'2025-06-17"
$ws.Range("F3").Style = "Normal"
